# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" sheet (before the "总计" summary sheet) with the
#    same layout/style as the existing quarterly sheets, populated with the
#    new quarter's fund-holding figures.
# 2) Update the "总计" (totals) sheet: add a new top data row for 2022-Q1
#    (pushing the existing quarters down) and update the running totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" quarterly sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($total)
$q1.Name = "2022-Q1"

# Reuse the layout/formatting of the previous quarter's sheet as a template.
$prev = $wb.Worksheets.Item("2021-Q4")
$prev.Range("B1:H2").Copy()
$q1.Range("B1").PasteSpecial(-4122)
$prev.Range("A2").Copy()
$q1.Range("A2").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Keep the textual columns as plain text (matching the source data, which
# preserves things like leading zeros in fund codes).
$q1.Range("B2:G2").NumberFormat = "@"

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "007280"
$q1.Range("C2").Value = "上投摩根日本精选股票（QDII）"
$q1.Range("D2").Value = "1.35"
$q1.Range("E2").Value = "88.71"
$q1.Range("F2").Value = "3.35"
$q1.Range("G2").Value = "0.0452"
$q1.Range("H2").Value = 5

# ---------------------------------------------------------------------
# 2) "总计" sheet: insert the new quarter at the top
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Rows.Item(2).Insert()
$ws.Range("B2:D2").ClearFormats()

$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q1"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0.05

$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3

# Restore the originally active sheet (unchanged by this edit).
$wb.Worksheets.Item("2021-Q2").Activate() | Out-Null
